$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text that can look like numbers (e.g. "211.81")
# Force those cells to stay text so Excel does not reinterpret them as numeric values,
# then reset the style so no stray number-format/quote-prefix style sticks to the cell.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.286.00"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.84%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.572.27"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.45%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.30%  "

$ws.Range("E6").Value = "  -0.68%  "

$ws.Range("E7").Value = "  +0.14%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "44.48"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.93%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "23.72"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.17%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.245"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.16%  "

$ws.Range("E11").Value = "  -1.09%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0894"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.50%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.798.37"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.30%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.569.79"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.64%  "

$ws.Range("E15").Value = "  -0.77%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.515"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.45%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "28.312.23"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.85%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "61.57"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.40%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "229.79"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.13%  "

$ws.Range("E20").Value = "  -0.22%  "

$ws.Range("E21").Value = "  -1.45%  "

$ws.Range("E22").Value = "  +0.07%  "

$ws.Range("E23").Value = "  +0.21%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.03"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.66%  "

$ws.Range("E25").Value = "  +0.45%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.55"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.07%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.93"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.82%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.35"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.95%  "

$ws.Range("E29").Value = "  -1.92%  "

$ws.Range("E30").Value = "  +0.07%  "

$ws.Range("E31").Value = "  +2.74%  "

$ws.Range("E32").Value = "  -3.48%  "

$ws.Range("E33").Value = "  -0.75%  "

$ws.Range("E34").Value = "  -1.93%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.378.10"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.47%  "

$ws.Range("E36").Value = "  +5.19%  "

$ws.Range("E38").Value = "  -0.06%  "

$ws.Range("E39").Value = "  +1.68%  "

$ws.Range("E40").Value = "  -1.82%  "

$ws.Range("E41").Value = "  -2.57%  "

$ws.Range("E42").Value = "  +0.02%  "

$ws.Range("E43").Value = "  +1.75%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.784"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.37%  "

$ws.Range("E46").Value = "  -4.19%  "

$ws.Range("E47").Value = "  -5.71%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "62.16"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.87%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.709.92"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.31%  "

$ws.Range("E50").Value = "  +0.88%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "85.47"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.52%  "
